# Tasks.xlsx — "Updated agenda and tasks after meeting"
#
# A new "Group meeting" task took place on 16.03 (Week 5), attended by
# Georgi, Mikaeil and Ilia, lasting 33 minutes. Record it as a new row
# right after the existing "Process report" row (row 24), in the same
# place the sheet already had a blank, pre-formatted row waiting for the
# next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 was a blank row formatted like the other task rows (e.g. row 3:
# A/B/C/D use the "task" look with borders/fill, D is numeric minutes).
# Copy that formatting down onto row 25 so the new entry matches the rest
# of the table exactly, then fill in the values.
$ws.Range("A3:D3").Copy() | Out-Null
$ws.Range("A25:D25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the "Week" label first so the shared-string table is populated in
# the same order as the source workbook (Week 5, then the task text,
# participants, and timestamp).
$ws.Range("E25").Value = "Week 5"
$ws.Range("A25").Value = "Group meeting"
$ws.Range("B25").Value = "Georgi, Mikaeil and Ilia"
$ws.Range("C25").Value = "16.03. || 11:00"
$ws.Range("D25").Value = 33

# Leave the selection where the editor ended up after entering the data.
$ws.Range("H1").Select() | Out-Null
